$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Registrar las horas consumidas en el Día 2 para la tarea (K6)
$ws.Range("K6").Value = 4

# Actualizar la celda seleccionada en la vista de la hoja
$ws.Range("E10").Select()
